$d = $word.ActiveDocument

# --- 1. Mark the five previously-unmarked drawing runs as "no proofing"
#        (adds <w:rPr><w:noProof/></w:rPr> in front of the <w:drawing>/
#        <w:lastRenderedPageBreak/> the same way Word does when it last
#        regenerated those runs). These are InlineShapes 7-11 (anchorIds
#        40A7828A, 40C2D0A8, 36BDCF84, 286C2720 and 3A1018CA).
for ($i = 7; $i -le 11; $i++) {
    $shape = $d.InlineShapes($i)
    $shape.Range.NoProofing = -1
}

# --- 2. Remove the stray _GoBack bookmark around
#        "CADA ENTIDADE TEM SUA CLASSE VALIDATION"
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# --- 3. Append two new paragraphs at the end of the document body.
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertAfter("C# TO JSON SITE PARA SERIALIZAR CLASSE PARA JSON")

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$end = $d.Content
$end.Collapse(0)
$end.InsertAfter("DESABILITAR FORMATAÇÃO E VALIDAÇÃO DE ERROS AUTOMÁTICO NO WEB API")
